$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current last (thick-bottom-border) row, which
# pushes that last row down from row 5 to row 6 and opens up a fresh row 5
# with the same column structure.
$ws.Rows("5:5").Insert()

# Copy the formatting (borders/styles) of row 4 (a regular data row) onto
# the newly inserted row 5, so A5/B5 end up with the same style as the
# other interior rows.
$ws.Range("A4:B4").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new people / tasks.
$ws.Range("A5").Value = "Ростислав Бердниченко"
$ws.Range("B5").Value = "Development"
$ws.Range("A6").Value = "Науменко Артем"
$ws.Range("B6").Value = "Writing documentation"

# Widen column A to fit the new, longer names (Excel stores column width
# padded by ~5/6 of a character vs. the ColumnWidth COM property, so this
# value serializes to width="22" in the saved XML).
$ws.Columns("A:A").ColumnWidth = 21.166666666666668

# Update the selected cell, matching the author's final cursor position.
[void]$ws.Range("I14").Select()
